$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item(1)
$ws.Range("H17").Value = 30885.445
$ws.Range("J17").Value = 30885.445
$ws.Range("L17").Value = 92656.33499999999
$ws.Range("N17").Value = -92992.33499999999
$ws.Range("H41").Value = 1554.1111
$ws.Range("I41").Value = 1947.8889
$ws.Range("J41").Value = 1160.3334
$ws.Range("K41").Value = 1947.8889
$ws.Range("L41").Value = 1160.3334
$ws.Range("M41").Value = -1507.8889
$ws.Range("N41").Value = -2040.3334
$ws.Range("H51").Value = 13999.9
$ws.Range("I51").Value = 4999.5
$ws.Range("J51").Value = 16250
$ws.Range("K51").Value = 4999.5
$ws.Range("L51").Value = 16250
$ws.Range("M51").Value = -4515.5
$ws.Range("N51").Value = -17218
$ws.Range("H64").Value = 5664.8335
$ws.Range("J64").Value = 6092.625
$ws.Range("L64").Value = 6092.625
$ws.Range("N64").Value = -6588.625
$ws.Range("H67").Value = 5664.8335
$ws.Range("J67").Value = 6092.625
$ws.Range("L67").Value = 6092.625
$ws.Range("N67").Value = -7808.625
$ws.Range("H88").Value = 2646.1853
$ws.Range("J88").Value = 3518.75
$ws.Range("L88").Value = 3518.75
$ws.Range("N88").Value = -4330.75
$ws.Range("H91").Value = 2646.1853
$ws.Range("J91").Value = 3518.75
$ws.Range("L91").Value = 3518.75
$ws.Range("N91").Value = -6326.75
$ws.Range("H98").Value = 1443.1111
$ws.Range("I98").Value = 1051.44
$ws.Range("K98").Value = 1051.44
$ws.Range("M98").Value = 446.5599999999999
$ws.Range("H112").Value = 5402.082
$ws.Range("J112").Value = 5422.1333
$ws.Range("L112").Value = 16266.3999
$ws.Range("N112").Value = -18482.3999
$ws.Range("H113").Value = 4046.5386
$ws.Range("I113").Value = 2976
$ws.Range("J113").Value = 5295.5
$ws.Range("K113").Value = 2976
$ws.Range("L113").Value = 5295.5
$ws.Range("M113").Value = 278
$ws.Range("N113").Value = -11803.5
$ws.Range("H115").Value = 2744.4443
$ws.Range("I115").Value = 671.4286
$ws.Range("K115").Value = 2014.2858
$ws.Range("M115").Value = -447.2857999999999
$ws.Range("H118").Value = 1522.8667
$ws.Range("I118").Value = 1545.6154
$ws.Range("K118").Value = 4636.8462
$ws.Range("M118").Value = -2979.8462
$ws.Range("H122").Value = 1443.1111
$ws.Range("I122").Value = 1051.44
$ws.Range("K122").Value = 3154.32
$ws.Range("M122").Value = -704.3200000000002
$ws.Range("H129").Value = 115840.11
$ws.Range("I129").Value = 982.2222
$ws.Range("K129").Value = 2946.6666
$ws.Range("M129").Value = 2053.3334
$ws.Range("H132").Value = 2000.5211
$ws.Range("I132").Value = 1991.8088
$ws.Range("K132").Value = 5975.4264
$ws.Range("M132").Value = -3445.4264
$ws.Range("H134").Value = 85861.69
$ws.Range("J134").Value = 85861.69
$ws.Range("L134").Value = 85861.69
$ws.Range("N134").Value = -96001.69
$ws.Range("H137").Value = 2093.3137
$ws.Range("I137").Value = 1972.4445
$ws.Range("J137").Value = 2229.2917
$ws.Range("K137").Value = 5917.333500000001
$ws.Range("L137").Value = 6687.875100000001
$ws.Range("M137").Value = -3367.333500000001
$ws.Range("N137").Value = -11787.8751
$ws.Range("H141").Value = 5682.3335
$ws.Range("I141").Value = 5682.3335
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 17047.0005
$ws.Range("L141").Value = 0
$ws.Range("M141").ClearContents()
$ws.Range("N141").Value = -11867.0005

# Sheet: ARM
$ws = $wb.Worksheets.Item(2)
$ws.Range("H2").Value = 1100.3055
$ws.Range("I2").Value = 1113.5588
$ws.Range("J2").Value = 875
$ws.Range("K2").Value = 1113.5588
$ws.Range("L2").Value = 875
$ws.Range("M2").Value = -1000.5588
$ws.Range("N2").Value = -1101
$ws.Range("H28").Value = 7427.857
$ws.Range("I28").Value = 6357.5
$ws.Range("K28").Value = 6357.5
$ws.Range("M28").Value = -6165.5
$ws.Range("H32").Value = 1172623.1
$ws.Range("I32").Value = 16967.684
$ws.Range("K32").Value = 16967.684
$ws.Range("M32").Value = -16680.684
$ws.Range("H45").Value = 3379.878
$ws.Range("I45").Value = 2712.7568
$ws.Range("J45").Value = 3928.4
$ws.Range("K45").Value = 2712.7568
$ws.Range("L45").Value = 3928.4
$ws.Range("M45").Value = -2335.7568
$ws.Range("N45").Value = -4682.4
$ws.Range("H74").Value = 4794.3184
$ws.Range("I74").Value = 3322.1177
$ws.Range("J74").Value = 9799.799999999999
$ws.Range("K74").Value = 3322.1177
$ws.Range("L74").Value = 9799.799999999999
$ws.Range("M74").Value = -2448.1177
$ws.Range("N74").Value = -11547.8
$ws.Range("H77").Value = 4794.3184
$ws.Range("I77").Value = 3322.1177
$ws.Range("J77").Value = 9799.799999999999
$ws.Range("K77").Value = 16610.5885
$ws.Range("L77").Value = 48999
$ws.Range("M77").Value = -12242.5885
$ws.Range("N77").Value = -57735
$ws.Range("H88").Value = 2707.5
$ws.Range("I88").Value = 2022.25
$ws.Range("J88").Value = 3164.3333
$ws.Range("K88").Value = 2022.25
$ws.Range("L88").Value = 3164.3333
$ws.Range("M88").Value = -1616.25
$ws.Range("N88").Value = -3976.3333
$ws.Range("H91").Value = 2707.5
$ws.Range("I91").Value = 2022.25
$ws.Range("J91").Value = 3164.3333
$ws.Range("K91").Value = 2022.25
$ws.Range("L91").Value = 3164.3333
$ws.Range("M91").Value = -618.25
$ws.Range("N91").Value = -5972.3333
$ws.Range("H98").Value = 24332.334
$ws.Range("J98").Value = 24332.334
$ws.Range("L98").Value = 24332.334
$ws.Range("N98").Value = -30322.334
$ws.Range("H99").Value = 7427.857
$ws.Range("I99").Value = 6357.5
$ws.Range("K99").Value = 6357.5
$ws.Range("M99").Value = -3362.5
$ws.Range("H103").Value = 98465.5
$ws.Range("J103").Value = 98465.5
$ws.Range("L103").Value = 98465.5
$ws.Range("N103").Value = -100809.5
$ws.Range("H110").Value = 1801.5186
$ws.Range("I110").Value = 1635.9048
$ws.Range("K110").Value = 1635.9048
$ws.Range("M110").Value = 409.0952
$ws.Range("H116").Value = 1100.3055
$ws.Range("I116").Value = 1113.5588
$ws.Range("J116").Value = 875
$ws.Range("K116").Value = 1113.5588
$ws.Range("L116").Value = 875
$ws.Range("M116").Value = 1180.4412
$ws.Range("N116").Value = -5463
$ws.Range("H122").Value = 6378.1177
$ws.Range("I122").Value = 8101.476
$ws.Range("J122").Value = 3594.2307
$ws.Range("K122").Value = 24304.428
$ws.Range("L122").Value = 10782.6921
$ws.Range("M122").Value = -21854.428
$ws.Range("N122").Value = -15682.6921
$ws.Range("H132").Value = 2307.5715
$ws.Range("I132").Value = 1712.8462
$ws.Range("K132").Value = 5138.5386
$ws.Range("M132").Value = -2608.5386
$ws.Range("H135").Value = 53144.59
$ws.Range("J135").Value = 53144.59
$ws.Range("L135").Value = 53144.59
$ws.Range("N135").Value = -63284.59
$ws.Range("H137").Value = 89998
$ws.Range("J137").Value = 89998
$ws.Range("L137").Value = 89998
$ws.Range("N137").Value = -100198
$ws.Range("H139").Value = 59084.43
$ws.Range("J139").Value = 59084.43
$ws.Range("L139").Value = 59084.43
$ws.Range("N139").Value = -69364.42999999999

# Sheet: BSM
$ws = $wb.Worksheets.Item(3)
$ws.Range("H2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("L2").ClearContents()
$ws.Range("N2").Value = 0
$ws.Range("H3").Value = 1100.3055
$ws.Range("I3").Value = 1113.5588
$ws.Range("J3").Value = 875
$ws.Range("K3").Value = 1113.5588
$ws.Range("L3").Value = 875
$ws.Range("M3").Value = -999.5588
$ws.Range("N3").Value = -1103
$ws.Range("H6").Value = 13849.5
$ws.Range("J6").Value = 13849.5
$ws.Range("L6").Value = 13849.5
$ws.Range("N6").Value = -14075.5
$ws.Range("H80").Value = 547.3333
$ws.Range("I80").Value = 222.33333
$ws.Range("J80").Value = 655.6667
$ws.Range("K80").Value = 222.33333
$ws.Range("L80").Value = 655.6667
$ws.Range("M80").Value = 775.6666700000001
$ws.Range("N80").Value = -2651.6667
$ws.Range("H83").Value = 547.3333
$ws.Range("I83").Value = 222.33333
$ws.Range("J83").Value = 655.6667
$ws.Range("K83").Value = 1111.66665
$ws.Range("L83").Value = 3278.3335
$ws.Range("M83").Value = 3880.33335
$ws.Range("N83").Value = -13262.3335
$ws.Range("H99").Value = 2534.7036
$ws.Range("I99").Value = 2439.8845
$ws.Range("K99").Value = 2439.8845
$ws.Range("M99").Value = -941.8845000000001
$ws.Range("H100").Value = 32085.777
$ws.Range("J100").Value = 32085.777
$ws.Range("L100").Value = 32085.777
$ws.Range("N100").Value = -34249.777
$ws.Range("H134").Value = 1602.52
$ws.Range("I134").Value = 1602.52
$ws.Range("K134").Value = 4807.559999999999
$ws.Range("M134").Value = -2272.559999999999

# Sheet: CRP
$ws = $wb.Worksheets.Item(4)
$ws.Range("H7").Value = 411.42856
$ws.Range("I7").Value = 658.75
$ws.Range("K7").Value = 658.75
$ws.Range("M7").Value = -545.75
$ws.Range("H22").Value = 1164.25
$ws.Range("I22").Value = 219
$ws.Range("K22").Value = 219
$ws.Range("M22").Value = 131
$ws.Range("H31").Value = 7868.9473
$ws.Range("I31").Value = 4561.7
$ws.Range("K31").Value = 4561.7
$ws.Range("M31").Value = -4266.7
$ws.Range("H34").Value = 7868.9473
$ws.Range("I34").Value = 4561.7
$ws.Range("K34").Value = 4561.7
$ws.Range("M34").Value = -4359.7
$ws.Range("H94").Value = 2650
$ws.Range("J94").Value = 2650
$ws.Range("L94").Value = 2650
$ws.Range("N94").Value = -3552
$ws.Range("I107").Value = 29412550
$ws.Range("J107").Value = 2005.4
$ws.Range("K107").Value = 29412550
$ws.Range("L107").Value = 2005.4
$ws.Range("M107").Value = -29410630
$ws.Range("N107").Value = -5845.4
$ws.Range("H115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("L115").ClearContents()
$ws.Range("N115").Value = 0
$ws.Range("H117").Value = 59999
$ws.Range("J117").Value = 59999
$ws.Range("L117").Value = 59999
$ws.Range("N117").Value = -69177
$ws.Range("H132").Value = 1072.2572
$ws.Range("I132").Value = 1007.0645
$ws.Range("J132").Value = 1577.5
$ws.Range("K132").Value = 3021.1935
$ws.Range("L132").Value = 4732.5
$ws.Range("M132").Value = -491.1934999999999
$ws.Range("N132").Value = -9792.5
$ws.Range("H134").Value = 943.2
$ws.Range("I134").Value = 825.8889
$ws.Range("K134").Value = 2477.6667
$ws.Range("M134").Value = 57.33329999999978
$ws.Range("H140").Value = 85038
$ws.Range("J140").Value = 85038
$ws.Range("L140").Value = 85038
$ws.Range("N140").Value = -95398
$ws.Range("H141").Value = 201371.69
$ws.Range("J141").Value = 232682.45
$ws.Range("L141").Value = 232682.45
$ws.Range("N141").Value = -243042.45

# Sheet: CUL
$ws = $wb.Worksheets.Item(5)
$ws.Range("H8").Value = 111617.75
$ws.Range("I8").Value = 111617.75
$ws.Range("K8").Value = 334853.25
$ws.Range("M8").Value = -334714.25
$ws.Range("H12").Value = 1000
$ws.Range("J12").Value = 1000
$ws.Range("L12").Value = 3000
$ws.Range("N12").Value = -3346
$ws.Range("H23").Value = 357.4
$ws.Range("J23").Value = 357.4
$ws.Range("L23").Value = 1072.2
$ws.Range("N23").Value = -1542.2
$ws.Range("H39").Value = 24999
$ws.Range("J39").Value = 24999
$ws.Range("L39").Value = 74997
$ws.Range("N39").Value = -75585
$ws.Range("H68").Value = 1000.2222
$ws.Range("J68").Value = 1000.2222
$ws.Range("L68").Value = 3000.6666
$ws.Range("N68").Value = -4622.6666
$ws.Range("H71").Value = 1000.2222
$ws.Range("J71").Value = 1000.2222
$ws.Range("L71").Value = 9001.9998
$ws.Range("N71").Value = -17113.9998
$ws.Range("H80").Value = 5554.5
$ws.Range("J80").Value = 5554.5
$ws.Range("L80").Value = 16663.5
$ws.Range("N80").Value = -18535.5
$ws.Range("H83").Value = 5554.5
$ws.Range("J83").Value = 5554.5
$ws.Range("L83").Value = 49990.5
$ws.Range("N83").Value = -59350.5
$ws.Range("H107").Value = 1140.909
$ws.Range("J107").Value = 1279.1428
$ws.Range("L107").Value = 3837.4284
$ws.Range("N107").Value = -7677.428400000001
$ws.Range("H112").Value = 11769.786
$ws.Range("I112").Value = 4557.6
$ws.Range("J112").Value = 15776.556
$ws.Range("K112").Value = 13672.8
$ws.Range("L112").Value = 47329.66800000001
$ws.Range("M112").Value = -12564.8
$ws.Range("N112").Value = -49545.66800000001
$ws.Range("H122").Value = 2367.0667
$ws.Range("I122").Value = 862.8333
$ws.Range("J122").Value = 2743.125
$ws.Range("K122").Value = 7765.4997
$ws.Range("L122").Value = 24688.125
$ws.Range("M122").Value = -5315.4997
$ws.Range("N122").Value = -29588.125
$ws.Range("H133").Value = 12732.25
$ws.Range("I133").Value = 8671.857
$ws.Range("J133").Value = 14918.615
$ws.Range("K133").Value = 26015.571
$ws.Range("L133").Value = 44755.845
$ws.Range("M133").Value = -20955.571
$ws.Range("N133").Value = -54875.845
$ws.Range("H134").Value = 4888.353
$ws.Range("I134").Value = 1650.2858
$ws.Range("J134").Value = 19999.334
$ws.Range("K134").Value = 4950.857400000001
$ws.Range("L134").Value = 59998.00199999999
$ws.Range("M134").Value = 119.1425999999992
$ws.Range("N134").Value = -70138.00199999999
$ws.Range("H137").Value = 4838.9375
$ws.Range("J137").Value = 5871.75
$ws.Range("L137").Value = 17615.25
$ws.Range("N137").Value = -27815.25

# Sheet: GSM
$ws = $wb.Worksheets.Item(6)
$ws.Range("H2").Value = 117.53846
$ws.Range("I2").Value = 151.28572
$ws.Range("J2").Value = 78.166664
$ws.Range("K2").Value = 151.28572
$ws.Range("L2").Value = 78.166664
$ws.Range("M2").Value = -38.28572
$ws.Range("N2").Value = -304.166664
$ws.Range("H17").Value = 3499.6667
$ws.Range("I17").Value = 7499
$ws.Range("K17").Value = 7499
$ws.Range("M17").Value = -7331
$ws.Range("H42").Value = 93931.2
$ws.Range("J42").Value = 93931.2
$ws.Range("L42").Value = 93931.2
$ws.Range("N42").Value = -94901.2
$ws.Range("H80").Value = 7798.727
$ws.Range("I80").Value = 7332.1665
$ws.Range("K80").Value = 7332.1665
$ws.Range("M80").Value = -6334.1665
$ws.Range("H83").Value = 7798.727
$ws.Range("I83").Value = 7332.1665
$ws.Range("K83").Value = 36660.8325
$ws.Range("M83").Value = -31668.8325
$ws.Range("H101").Value = 82049.586
$ws.Range("J101").Value = 82049.586
$ws.Range("L101").Value = 82049.586
$ws.Range("N101").Value = -88539.586
$ws.Range("H102").Value = 1843.3334
$ws.Range("I102").Value = 1865.75
$ws.Range("K102").Value = 1865.75
$ws.Range("M102").Value = -243.75
$ws.Range("H107").Value = 1360
$ws.Range("I107").Value = 1360
$ws.Range("K107").Value = 1360
$ws.Range("M107").Value = 560
$ws.Range("H113").Value = 4828.8125
$ws.Range("I113").Value = 6004.3335
$ws.Range("J113").Value = 1302.25
$ws.Range("K113").Value = 6004.3335
$ws.Range("L113").Value = 1302.25
$ws.Range("M113").Value = -3834.3335
$ws.Range("N113").Value = -5642.25
$ws.Range("H115").Value = 93931.2
$ws.Range("J115").Value = 93931.2
$ws.Range("L115").Value = 93931.2
$ws.Range("N115").Value = -96281.2
$ws.Range("H122").Value = 2611.7222
$ws.Range("I122").Value = 1993.1538
$ws.Range("J122").Value = 4220
$ws.Range("K122").Value = 5979.4614
$ws.Range("L122").Value = 12660
$ws.Range("M122").Value = -3529.4614
$ws.Range("N122").Value = -17560
$ws.Range("H132").Value = 3728.717
$ws.Range("I132").Value = 3378.3
$ws.Range("K132").Value = 10134.9
$ws.Range("M132").Value = -7604.900000000001
$ws.Range("H135").Value = 71448.234
$ws.Range("J135").Value = 71448.234
$ws.Range("L135").Value = 71448.234
$ws.Range("N135").Value = -81588.234
$ws.Range("H136").Value = 36447.54
$ws.Range("J136").Value = 36447.54
$ws.Range("L136").Value = 109342.62
$ws.Range("N136").Value = -114442.62

# Sheet: LTW
$ws = $wb.Worksheets.Item(7)
$ws.Range("H16").Value = 1449.1538
$ws.Range("I16").Value = 1449.1538
$ws.Range("K16").Value = 1449.1538
$ws.Range("M16").Value = -1279.1538
$ws.Range("H22").Value = 1759.7222
$ws.Range("I22").Value = 1221.3636
$ws.Range("J22").Value = 2605.7144
$ws.Range("K22").Value = 1221.3636
$ws.Range("L22").Value = 2605.7144
$ws.Range("M22").Value = -926.3635999999999
$ws.Range("N22").Value = -3195.7144
$ws.Range("H27").Value = 1759.7222
$ws.Range("I27").Value = 1221.3636
$ws.Range("J27").Value = 2605.7144
$ws.Range("K27").Value = 1221.3636
$ws.Range("L27").Value = 2605.7144
$ws.Range("M27").Value = -1114.3636
$ws.Range("N27").Value = -2819.7144
$ws.Range("H55").Value = 392.26666
$ws.Range("I55").Value = 276.41666
$ws.Range("J55").Value = 855.6667
$ws.Range("K55").Value = 276.41666
$ws.Range("L55").Value = 855.6667
$ws.Range("M55").Value = -103.41666
$ws.Range("N55").Value = -1201.6667
$ws.Range("H61").Value = 1575.6072
$ws.Range("I61").Value = 1442.5
$ws.Range("K61").Value = 1442.5
$ws.Range("M61").Value = -1240.5
$ws.Range("H82").Value = 5166.7144
$ws.Range("J82").Value = 2994.5
$ws.Range("L82").Value = 2994.5
$ws.Range("N82").Value = -3716.5
$ws.Range("H85").Value = 5166.7144
$ws.Range("J85").Value = 2994.5
$ws.Range("L85").Value = 2994.5
$ws.Range("N85").Value = -5490.5
$ws.Range("H100").Value = 3124.8572
$ws.Range("I100").Value = 3124.8572
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 3124.8572
$ws.Range("L100").Value = 0
$ws.Range("M100").ClearContents()
$ws.Range("N100").Value = -2583.8572
$ws.Range("H101").Value = 52589.75
$ws.Range("J101").Value = 52589.75
$ws.Range("L101").Value = 52589.75
$ws.Range("N101").Value = -59079.75
$ws.Range("H113").Value = 1575.6072
$ws.Range("I113").Value = 1442.5
$ws.Range("K113").Value = 1442.5
$ws.Range("M113").Value = 727.5
$ws.Range("H117").Value = 0
$ws.Range("I117").Value = 0
$ws.Range("K117").Value = 0
$ws.Range("M117").ClearContents()
$ws.Range("H122").Value = 3057.077
$ws.Range("I122").Value = 2758.3
$ws.Range("K122").Value = 8274.900000000001
$ws.Range("M122").Value = -5824.900000000001
$ws.Range("H132").Value = 6847.604
$ws.Range("I132").Value = 7507.75
$ws.Range("J132").Value = 6108.24
$ws.Range("K132").Value = 22523.25
$ws.Range("L132").Value = 18324.72
$ws.Range("M132").Value = -19993.25
$ws.Range("N132").Value = -23384.72
$ws.Range("H133").Value = 86665.59
$ws.Range("J133").Value = 86665.59
$ws.Range("L133").Value = 86665.59
$ws.Range("N133").Value = -91725.59
$ws.Range("H134").Value = 87373.125
$ws.Range("J134").Value = 87373.125
$ws.Range("L134").Value = 87373.125
$ws.Range("N134").Value = -97513.125
$ws.Range("H136").Value = 4060.9333
$ws.Range("I136").Value = 3882.5557
$ws.Range("K136").Value = 11647.6671
$ws.Range("M136").Value = -9097.667099999999
$ws.Range("H139").Value = 89438.766
$ws.Range("J139").Value = 89423.27
$ws.Range("L139").Value = 89423.27
$ws.Range("N139").Value = -99703.27
$ws.Range("H141").Value = 89925.7
$ws.Range("J141").Value = 89925.7
$ws.Range("L141").Value = 89925.7
$ws.Range("N141").Value = -100285.7

# Sheet: WVR
$ws = $wb.Worksheets.Item(8)
$ws.Range("H24").Value = 13000
$ws.Range("J24").Value = 13000
$ws.Range("L24").Value = 13000
$ws.Range("N24").Value = -13460
$ws.Range("H44").Value = 20000
$ws.Range("J44").Value = 20000
$ws.Range("L44").Value = 20000
$ws.Range("N44").Value = -21108
$ws.Range("H62").Value = 8998.125
$ws.Range("I62").Value = 7997.4
$ws.Range("K62").Value = 7997.4
$ws.Range("M62").Value = -7373.4
$ws.Range("H65").Value = 8998.125
$ws.Range("I65").Value = 7997.4
$ws.Range("K65").Value = 39987
$ws.Range("M65").Value = -36867
$ws.Range("H104").Value = 10625.167
$ws.Range("J104").Value = 10625.167
$ws.Range("L104").Value = 10625.167
$ws.Range("N104").Value = -17613.167
$ws.Range("H113").Value = 590.3125
$ws.Range("I113").Value = 674.8333
$ws.Range("J113").Value = 336.75
$ws.Range("K113").Value = 2024.4999
$ws.Range("L113").Value = 1010.25
$ws.Range("M113").Value = 145.5001
$ws.Range("N113").Value = -5350.25
$ws.Range("H118").Value = 0
$ws.Range("I118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("K118").Value = 0
$ws.Range("L118").ClearContents()
$ws.Range("M118").ClearContents()
$ws.Range("N118").Value = 0
$ws.Range("H122").Value = 3232.963
$ws.Range("I122").Value = 2425.4348
$ws.Range("K122").Value = 7276.3044
$ws.Range("M122").Value = -4826.3044
$ws.Range("H129").Value = 99429
$ws.Range("J129").Value = 99429
$ws.Range("L129").Value = 99429
$ws.Range("N129").Value = -109429
$ws.Range("H132").Value = 3555.8857
$ws.Range("I132").Value = 3238.1724
$ws.Range("K132").Value = 9714.5172
$ws.Range("M132").Value = -7184.5172
$ws.Range("H136").Value = 11066.421
$ws.Range("I136").Value = 11829.078
$ws.Range("K136").Value = 35487.234
$ws.Range("M136").Value = -32937.234
